# This script updates the LR-pairs worksheet (Ccl21b-Ccr10) with refreshed TPM-based
# NATMI statistics, adding a new "Inflammatory-Mac" sending cluster (rows 8-10 data is
# now Inflammatory-Mac instead of MuSCs) and appending 3 new rows (11-13) for the MuSCs
# sending cluster against each target cluster (ECs, FAPs, MuSCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07908133333333334
$ws.Range("H2").Value = 0.237244
$ws.Range("I2").Value = 0.1380838080781507
$ws.Range("J2").Value = 0.1380838080781507
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.478362
$ws.Range("N2").Value = 1.435086
$ws.Range("O2").Value = 0.6533545125880439
$ws.Range("P2").Value = 0.6533545125880439
$ws.Range("Q2").Value = 0.037829504776
$ws.Range("R2").Value = 0.3404655429840001
$ws.Range("S2").Value = 0.09021767912320117
$ws.Range("T2").Value = 0.09021767912320115

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07908133333333334
$ws.Range("H3").Value = 0.237244
$ws.Range("I3").Value = 0.1380838080781507
$ws.Range("J3").Value = 0.1380838080781507
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.155562
$ws.Range("N3").Value = 0.466686
$ws.Range("O3").Value = 0.2124690813384451
$ws.Range("P3").Value = 0.2124690813384451
$ws.Range("Q3").Value = 0.012302050376
$ws.Range("R3").Value = 0.110718453384
$ws.Range("S3").Value = 0.02933853985007886
$ws.Range("T3").Value = 0.02933853985007885

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07908133333333334
$ws.Range("H4").Value = 0.237244
$ws.Range("I4").Value = 0.1380838080781507
$ws.Range("J4").Value = 0.1380838080781507
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.098239
$ws.Range("N4").Value = 0.294717
$ws.Range("O4").Value = 0.134176406073511
$ws.Range("P4").Value = 0.1341764060735109
$ws.Range("Q4").Value = 0.007768871105333334
$ws.Range("R4").Value = 0.06991983994800001
$ws.Range("S4").Value = 0.01852758910487071
$ws.Range("T4").Value = 0.0185275891048707

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07720866666666666
$ws.Range("H5").Value = 0.231626
$ws.Range("I5").Value = 0.1348139473702591
$ws.Range("J5").Value = 0.134813947370259
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.478362
$ws.Range("N5").Value = 1.435086
$ws.Range("O5").Value = 0.6533545125880439
$ws.Range("P5").Value = 0.6533545125880439
$ws.Range("Q5").Value = 0.036933692204
$ws.Range("R5").Value = 0.332403229836
$ws.Range("S5").Value = 0.08808130087416581
$ws.Range("T5").Value = 0.0880813008741658

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ccl21b"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.07720866666666666
$ws.Range("H6").Value = 0.231626
$ws.Range("I6").Value = 0.1348139473702591
$ws.Range("J6").Value = 0.134813947370259
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.155562
$ws.Range("N6").Value = 0.466686
$ws.Range("O6").Value = 0.2124690813384451
$ws.Range("P6").Value = 0.2124690813384451
$ws.Range("Q6").Value = 0.012010734604
$ws.Range("R6").Value = 0.108096611436
$ws.Range("S6").Value = 0.02864379554936844
$ws.Range("T6").Value = 0.02864379554936843

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ccl21b"
$ws.Range("C7").Value = "Ccr10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.07720866666666666
$ws.Range("H7").Value = 0.231626
$ws.Range("I7").Value = 0.1348139473702591
$ws.Range("J7").Value = 0.134813947370259
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.098239
$ws.Range("N7").Value = 0.294717
$ws.Range("O7").Value = 0.134176406073511
$ws.Range("P7").Value = 0.1341764060735109
$ws.Range("Q7").Value = 0.007584902204666667
$ws.Range("R7").Value = 0.068264119842
$ws.Range("S7").Value = 0.01808885094672482
$ws.Range("T7").Value = 0.01808885094672481

# Row 8: Inflammatory-Mac -> ECs
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Ccl21b"
$ws.Range("C8").Value = "Ccr10"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.03917266666666667
$ws.Range("H8").Value = 0.117518
$ws.Range("I8").Value = 0.06839933974190335
$ws.Range("J8").Value = 0.06839933974190333
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.478362
$ws.Range("N8").Value = 1.435086
$ws.Range("O8").Value = 0.6533545125880439
$ws.Range("P8").Value = 0.6533545125880439
$ws.Range("Q8").Value = 0.018738715172
$ws.Range("R8").Value = 0.168648436548
$ws.Range("S8").Value = 0.04468901727841528
$ws.Range("T8").Value = 0.04468901727841527

# Row 9: Inflammatory-Mac -> FAPs
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Ccl21b"
$ws.Range("C9").Value = "Ccr10"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.03917266666666667
$ws.Range("H9").Value = 0.117518
$ws.Range("I9").Value = 0.06839933974190335
$ws.Range("J9").Value = 0.06839933974190333
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.155562
$ws.Range("N9").Value = 0.466686
$ws.Range("O9").Value = 0.2124690813384451
$ws.Range("P9").Value = 0.2124690813384451
$ws.Range("Q9").Value = 0.006093778372
$ws.Range("R9").Value = 0.054844005348
$ws.Range("S9").Value = 0.01453274487911841
$ws.Range("T9").Value = 0.0145327448791184

# Row 10: Inflammatory-Mac -> MuSCs
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Ccl21b"
$ws.Range("C10").Value = "Ccr10"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.03917266666666667
$ws.Range("H10").Value = 0.117518
$ws.Range("I10").Value = 0.06839933974190335
$ws.Range("J10").Value = 0.06839933974190333
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.098239
$ws.Range("N10").Value = 0.294717
$ws.Range("O10").Value = 0.134176406073511
$ws.Range("P10").Value = 0.1341764060735109
$ws.Range("Q10").Value = 0.003848283600666667
$ws.Range("R10").Value = 0.034634552406
$ws.Range("S10").Value = 0.009177577584369661
$ws.Range("T10").Value = 0.009177577584369658

# Row 11: MuSCs -> ECs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Ccl21b"
$ws.Range("C11").Value = "Ccr10"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3772426666666667
$ws.Range("H11").Value = 1.131728
$ws.Range("I11").Value = 0.658702904809687
$ws.Range("J11").Value = 0.6587029048096869
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.478362
$ws.Range("N11").Value = 1.435086
$ws.Range("O11").Value = 0.6533545125880439
$ws.Range("P11").Value = 0.6533545125880439
$ws.Range("Q11").Value = 0.180458556512
$ws.Range("R11").Value = 1.624127008608
$ws.Range("S11").Value = 0.4303665153122617
$ws.Range("T11").Value = 0.4303665153122616

# Row 12: MuSCs -> FAPs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ccl21b"
$ws.Range("C12").Value = "Ccr10"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3772426666666667
$ws.Range("H12").Value = 1.131728
$ws.Range("I12").Value = 0.658702904809687
$ws.Range("J12").Value = 0.6587029048096869
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.155562
$ws.Range("N12").Value = 0.466686
$ws.Range("O12").Value = 0.2124690813384451
$ws.Range("P12").Value = 0.2124690813384451
$ws.Range("Q12").Value = 0.05868462371200001
$ws.Range("R12").Value = 0.528161613408
$ws.Range("S12").Value = 0.1399540010598795
$ws.Range("T12").Value = 0.1399540010598794

# Row 13: MuSCs -> MuSCs
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ccl21b"
$ws.Range("C13").Value = "Ccr10"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3772426666666667
$ws.Range("H13").Value = 1.131728
$ws.Range("I13").Value = 0.658702904809687
$ws.Range("J13").Value = 0.6587029048096869
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.098239
$ws.Range("N13").Value = 0.294717
$ws.Range("O13").Value = 0.134176406073511
$ws.Range("P13").Value = 0.1341764060735109
$ws.Range("Q13").Value = 0.03705994233066667
$ws.Range("R13").Value = 0.333539480976
$ws.Range("S13").Value = 0.0883823884375458
$ws.Range("T13").Value = 0.08838238843754577
